$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fix the date (timesheet was re-dated) ---
$ws.Range("A2").Value = 41807

# --- Row 4: new time-sheet entry (evening session) ---
$ws.Range("B4").Value = 0.90972222222222221
$ws.Range("C4").Value = 0.9375
# Match the existing time formatting used by the other Start/End Time cells
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
# D4 already carries the shared HOUR()/MINUTE() formula from the template and
# recalculates automatically once B4/C4 are populated.
$ws.Range("E4").Value = "Can now write a basic pdf"

# --- Row 6: new day of entries ---
$ws.Range("A6").Value = 41808
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 41808

$ws.Range("B6").Value = 0.85416666666666663
$ws.Range("C6").Value = 0.88888888888888884
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
# D6 already carries the shared formula from the template.
$ws.Range("E6").Value = "Can no iterate through the individual sheet and print to pdf... Understand how highlighting works"

# --- Cosmetic: the "Hours" column (D) formula cells were styled with the
# Segoe UI font; bring them in line with the sheet's normal Calibri font
# (keeping the same dark-gray color) ---
$ws.Range("D2:D21").Font.Name = "Calibri"

# --- Selection moved to E9 as the sheet was last edited there ---
$ws.Range("E9").Select()
